$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "(K)" labeled AI column (I) with average-excluding-self formulas (H) ---
$ws.Range("H2").Formula = "=SUM(C2:F2)/4"
$ws.Range("I2").Value = "AI1(K)"

$ws.Range("H3").Formula = "=SUM(B3,D3:F3)/4"
$ws.Range("I3").Value = "AI2(K)"

$ws.Range("H4").Formula = "=SUM(B4:C4,E4:F4)/4"
$ws.Range("I4").Value = "AI3(K)"

$ws.Range("H5").Formula = "=SUM(B5:D5,F5)/4"
$ws.Range("I5").Value = "AI4(K)"

$ws.Range("H6").Formula = "=SUM(B6:E6)/4"
$ws.Range("I6").Value = "AI5(K)"

# --- Row 8: remaining-budget formulas ---
$ws.Range("B8").Formula = "=10000-SUM(B3:B6)/4"
$ws.Range("C8").Formula = "=10000-SUM(C2,C4:C6)/4"
$ws.Range("D8").Formula = "=10000-SUM(D2:D3,D5:D6)/4"
$ws.Range("E8").Formula = "=10000-SUM(E2:E4,E6)/4"
$ws.Range("F8").Formula = "=10000-SUM(F2:F5)/4"

# --- Row 9: repeated headers ---
$ws.Range("B9").Value = "AI1"
$ws.Range("C9").Value = "AI2"
$ws.Range("D9").Value = "AI3"
$ws.Range("E9").Value = "AI4"
$ws.Range("F9").Value = "AI5"

# --- Ranking table (B13:C22), sorted descending by value ---
$ws.Range("B13").Value = "AI5(K)"
$ws.Range("C13").Value = 9267.5

$ws.Range("B14").Value = "AI5"
$ws.Range("C14").Value = 8916

$ws.Range("B15").Value = "AI2(K)"
$ws.Range("C15").Value = 6119.75

$ws.Range("B16").Value = "AI4(K)"
$ws.Range("C16").Value = 5533.75

$ws.Range("B17").Value = "AI2"
$ws.Range("C17").Value = 5311.25

$ws.Range("B18").Value = "AI4"
$ws.Range("C18").Value = 4678.5

$ws.Range("B19").Value = "AI3(K)"
$ws.Range("C19").Value = 3545.75

$ws.Range("B20").Value = "AI3"
$ws.Range("C20").Value = 2816.25

$ws.Range("B21").Value = "AI1(K)"
$ws.Range("C21").Value = 2184

$ws.Range("B22").Value = "AI1"
$ws.Range("C22").Value = 1627.25

# --- Sort state metadata on the ranking table ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C13:C22"), 0, 2)
$ws.Sort.SetRange($ws.Range("B13:C22"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# --- View state ---
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B13:C22").Select() | Out-Null

$wb.Save() | Out-Null
